# Flight Manifest.xlsx update: "Update material and add step 2"
# - Column M ("Baggage") switches from a text YES/NO indicator to a
#   numeric integer count, formatted with number format "0".
# - The now-unused "NO"/"YES" shared strings disappear automatically
#   once no cell references them anymore.
# - Selection moves to M6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the header cell (bold font) the integer number format too, so a
# distinct style (bold + numFmt "0") is produced for M1.
$ws.Range("M1").NumberFormat = "0"

# Data cells: integer number format + new baggage counts.
$ws.Range("M2:M6").NumberFormat = "0"
$ws.Range("M2").Value = 0
$ws.Range("M3").Value = 2
$ws.Range("M4").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("M6").Value = 2

# Move the active selection to M6, matching the saved state of the file.
$ws.Range("M6").Select()
